$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "41.651.80"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "2.476.84"
$ws.Range("E3").Value = "  +0.82%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "319.25"
$ws.Range("E5").Value = "  +1.44%  "
Set-TextValue $ws.Range("D6") "92.36"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("E7").Value = "  +0.98%  "
$ws.Range("E8").Value = "  +0.05%  "
Set-TextValue $ws.Range("D9") "0.512"
$ws.Range("E9").Value = "  +0.82%  "
Set-TextValue $ws.Range("D10") "0.0866"
$ws.Range("E10").Value = "  +8.78%  "
Set-TextValue $ws.Range("D11") "33.09"
$ws.Range("E11").Value = "  +2.10%  "
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "2.858.77"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("E14").Value = "  +0.73%  "
Set-TextValue $ws.Range("D15") "15.54"
$ws.Range("E15").Value = "  -1.65%  "
$ws.Range("D16").Value = "2.478.01"
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("E17").Value = "  +2.77%  "
$ws.Range("D18").Value = "41.603.83"
Set-TextValue $ws.Range("D19") "6.45"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").Value = "0.0₃0944"
$ws.Range("E20").Value = "  +1.01%  "
Set-TextValue $ws.Range("D21") "70.71"
$ws.Range("E21").Value = "  -0.04%  "
Set-TextValue $ws.Range("D22") "11.27"
$ws.Range("E22").Value = "  -0.53%  "
Set-TextValue $ws.Range("D23") "239.89"
$ws.Range("E23").Value = "  +1.13%  "
$ws.Range("E24").Value = "  +1.61%  "
$ws.Range("E25").Value = "  +2.76%  "
$ws.Range("E26").Value = "  +0.01%  "
Set-TextValue $ws.Range("D27") "24.94"
$ws.Range("E27").Value = "  +2.88%  "
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("E30").Value = "  +4.15%  "
Set-TextValue $ws.Range("D31") "156.89"
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("E32").Value = "  +0.61%  "
Set-TextValue $ws.Range("D34") "0.0765"
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("E35").Value = "  -0.57%  "
Set-TextValue $ws.Range("D36") "17.20"
$ws.Range("E36").Value = "  -1.07%  "
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D37") "0.116"
$ws.Range("E37").Value = "  +1.83%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D38") "1.84"
$ws.Range("E38").Value = "  +3.61%  "
Set-TextValue $ws.Range("D39") "2.88"
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("E40").Value = "  +2.03%  "
Set-TextValue $ws.Range("D41") "4.01"
$ws.Range("E41").Value = "  +2.15%  "
$ws.Range("E42").Value = "  +3.05%  "
$ws.Range("D43").Value = "1.992.55"
$ws.Range("E43").Value = "  +0.99%  "
$ws.Range("E44").Value = "  +1.05%  "
Set-TextValue $ws.Range("D45") "18.88"
$ws.Range("E45").Value = "  +2.00%  "
Set-TextValue $ws.Range("D46") "2.97"
$ws.Range("E46").Value = "  +2.63%  "
Set-TextValue $ws.Range("D47") "9.46"
$ws.Range("E47").Value = "  +5.92%  "
$ws.Range("D48").Value = "2.714.39"
$ws.Range("E48").Value = "  +0.65%  "
Set-TextValue $ws.Range("D49") "98.00"
$ws.Range("E49").Value = "  +1.78%  "
Set-TextValue $ws.Range("D50") "75.45"
$ws.Range("E50").Value = "  +5.30%  "
Set-TextValue $ws.Range("D51") "67.14"
